$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Component Analysis filtering: a new leading error value (Q0) is inserted
# into column B for each data row, shifting the existing Q0..Q8 values one
# column to the right (B->C, C->D, ... J->K) and dropping the old Q9 value
# (old column K) that no longer fits within the Q0..Q9 window.

# Row 2
$ws.Range("B2").Value = -0.8260246319213993
$ws.Range("C2").Value = -1.981077477814098
$ws.Range("D2").Value = -0.4414436339245075
$ws.Range("E2").Value = -0.7110977538981412
$ws.Range("F2").Value = -0.159566049999028
$ws.Range("G2").Value = -0.3854928109118805
$ws.Range("H2").Value = -0.2776770955153309
$ws.Range("I2").Value = -0.3740767161796326
$ws.Range("J2").Value = 0.1293341692733759
$ws.Range("K2").Value = 0.3878561111968251

# Row 3
$ws.Range("B3").Value = -0.1405287498260583
$ws.Range("C3").Value = -0.4101828697996921
$ws.Range("D3").Value = 0.1413488340994211
$ws.Range("E3").Value = -0.08457792681343129
$ws.Range("F3").Value = 0.0232377885831183
$ws.Range("G3").Value = -0.0731618320811834
$ws.Range("H3").Value = 0.4302490533718251
$ws.Range("I3").Value = 0.6887709952952742
$ws.Range("J3").Value = -1.03521158968304
$ws.Range("K3").Value = -0.3815723733203383

# Row 4
$ws.Range("B4").Value = 0.7815531178611421
$ws.Range("C4").Value = 0.5556263569482897
$ws.Range("D4").Value = 0.6634420723448393
$ws.Range("E4").Value = 0.5670424516805376
$ws.Range("F4").Value = 1.070453337133546
$ws.Range("G4").Value = 1.328975279056995
$ws.Range("H4").Value = -0.3950073059213186
$ws.Range("I4").Value = 0.2586319104413827
$ws.Range("J4").Value = 0.6221758453498138
$ws.Range("K4").Value = -0.2785999995298156

# Row 5
$ws.Range("B5").Value = 0.59085446968239
$ws.Range("C5").Value = 0.4944548490180882
$ws.Range("D5").Value = 0.9978657344710967
$ws.Range("E5").Value = 1.256387676394546
$ws.Range("F5").Value = -0.467594908583768
$ws.Range("G5").Value = 0.1860443077789333
$ws.Range("H5").Value = 0.5495882426873645
$ws.Range("I5").Value = -0.351187602192265
$ws.Range("J5").Value = 0.7191879547613879
$ws.Range("K5").Value = 0.142128040744926

# Row 6
$ws.Range("B6").Value = 0.9421746901075616
$ws.Range("C6").Value = 1.200696632031011
$ws.Range("D6").Value = -0.5232859529473031
$ws.Range("E6").Value = 0.1303532634153982
$ws.Range("F6").Value = 0.4938971983238294
$ws.Range("G6").Value = -0.4068786465558001
$ws.Range("H6").Value = 0.6634969103978527
$ws.Range("I6").Value = 0.08643699638139091
$ws.Range("J6").Value = -0.1191689877116253
$ws.Range("K6").Value = 0.8599759940094719

# Row 7
$ws.Range("B7").Value = -0.8773451131212686
$ws.Range("C7").Value = -0.2237058967585673
$ws.Range("D7").Value = 0.1398380381498639
$ws.Range("E7").Value = -0.7609378067297656
$ws.Range("F7").Value = 0.3094377502238873
$ws.Range("G7").Value = -0.2676221637925746
$ws.Range("H7").Value = -0.4732281478855908
$ws.Range("I7").Value = 0.5059168338355065
$ws.Range("J7").Value = -0.2304295095018555
$ws.Range("K7").Value = -0.2483993241560489

# Row 8
$ws.Range("B8").Value = 0.04951982135029001
$ws.Range("C8").Value = -0.8512560235293395
$ws.Range("D8").Value = 0.2191195334243133
$ws.Range("E8").Value = -0.3579403805921485
$ws.Range("F8").Value = -0.5635463646851647
$ws.Range("G8").Value = 0.4155986170359325
$ws.Range("H8").Value = -0.3207477263014294
$ws.Range("I8").Value = -0.3387175409556228
$ws.Range("J8").Value = -0.315312675415286
$ws.Range("K8").Value = -0.6804871111820268

# Row 9
$ws.Range("B9").Value = 0.2874422675287736
$ws.Range("C9").Value = -0.2896176464876882
$ws.Range("D9").Value = -0.4952236305807045
$ws.Range("E9").Value = 0.4839213511403928
$ws.Range("F9").Value = -0.2524249921969692
$ws.Range("G9").Value = -0.2703948068511625
$ws.Range("H9").Value = -0.2469899413108257
$ws.Range("I9").Value = -0.6121643770775664
$ws.Range("J9").Value = 0.04365514009338567
$ws.Range("K9").Value = -0.08721574084803801

# Row 10
$ws.Range("B10").Value = -0.5938852086106097
$ws.Range("C10").Value = 0.3852597731104875
$ws.Range("D10").Value = -0.3510865702268745
$ws.Range("E10").Value = -0.3690563848810678
$ws.Range("F10").Value = -0.345651519340731
$ws.Range("G10").Value = -0.7108259551074718
$ws.Range("H10").Value = -0.0550064379365196
$ws.Range("I10").Value = -0.1858773188779433
$ws.Range("J10").Value = -0.2953951178020652
$ws.Range("K10").Value = -1.125172904869618

# Row 11
$ws.Range("B11").Value = -0.2740769919852016
$ws.Range("C11").Value = -0.292046806639395
$ws.Range("D11").Value = -0.2686419410990581
$ws.Range("E11").Value = -0.6338163768657989
$ws.Range("F11").Value = 0.02200314030515318
$ws.Range("G11").Value = -0.1088677406362705
$ws.Range("H11").Value = -0.2183855395603924
$ws.Range("I11").Value = -1.048163326627946
$ws.Range("J11").Value = -0.4676375374592979
$ws.Range("K11").Value = -0.5510236388890779

# Row 12
$ws.Range("B12").Value = -0.09712588908289738
$ws.Range("C12").Value = -0.4623003248496382
$ws.Range("D12").Value = 0.193519192321314
$ws.Range("E12").Value = 0.0626483113798903
$ws.Range("F12").Value = -0.04686948754423159
$ws.Range("G12").Value = -0.8766472746117848
$ws.Range("H12").Value = -0.2961214854431372
$ws.Range("I12").Value = -0.3795075868729171
$ws.Range("J12").Value = -0.5974174644126757
$ws.Range("K12").Value = 0.9663603392234945

# Row 13
$ws.Range("B13").Value = 0.1560540901775642
$ws.Range("C13").Value = 0.02518320923614054
$ws.Range("D13").Value = -0.08433458968798135
$ws.Range("E13").Value = -0.9141123767555346
$ws.Range("F13").Value = -0.3335865875868869
$ws.Range("G13").Value = -0.4169726890166669
$ws.Range("H13").Value = -0.6348825665564255
$ws.Range("I13").Value = 0.9288952370797448
$ws.Range("J13").Value = -0.4070884690023188
$ws.Range("K13").Value = 1.033525947527809

# Row 14
$ws.Range("B14").Value = -0.106237102096024
$ws.Range("C14").Value = -0.9360148891635773
$ws.Range("D14").Value = -0.3554890999949296
$ws.Range("E14").Value = -0.4388752014247095
$ws.Range("F14").Value = -0.6567850789644682
$ws.Range("G14").Value = 0.9069927246717022
$ws.Range("H14").Value = -0.4289909814103615
$ws.Range("I14").Value = 1.011623435119767
$ws.Range("J14").Value = -0.0888206791410559
$ws.Range("K14").Value = -0.2990950773107258

# Row 15
$ws.Range("B15").Value = -0.2673832081489602
$ws.Range("C15").Value = -0.3507693095787401
$ws.Range("D15").Value = -0.5686791871184989
$ws.Range("E15").Value = 0.9950986165176716
$ws.Range("F15").Value = -0.3408850895643921
$ws.Range("G15").Value = 1.099729326965736
$ws.Range("H15").Value = -0.0007147872950865053
$ws.Range("I15").Value = -0.2109891854647564
$ws.Range("J15").Value = 0.387555017918974
$ws.Range("K15").Value = 0.06188637714375961

# Row 16
$ws.Range("B16").Value = -0.4749660705519546
$ws.Range("C16").Value = 1.088811733084216
$ws.Range("D16").Value = -0.2471719729978479
$ws.Range("E16").Value = 1.19344244353228
$ws.Range("F16").Value = 0.09299832927145768
$ws.Range("G16").Value = -0.1172760688982122
$ws.Range("H16").Value = 0.4812681344855182
$ws.Range("I16").Value = 0.1555994937103038
$ws.Range("J16").Value = 0.8727792950859075
$ws.Range("K16").Value = 2.568651570723608

# Row 17
$ws.Range("B17").Value = 1.097015279354651
$ws.Range("C17").Value = -0.2389684267274125
$ws.Range("D17").Value = 1.201645989802716
$ws.Range("E17").Value = 0.1012018755418931
$ws.Range("F17").Value = -0.1090725226277768
$ws.Range("G17").Value = 0.4894716807559536
$ws.Range("H17").Value = 0.1638030399807392
$ws.Range("I17").Value = 0.8809828413563429
$ws.Range("J17").Value = 2.576855116994044
$ws.Range("K17").Value = 9.421101911918621

# Row 18
$ws.Range("B18").Value = -0.2415579873788807
$ws.Range("C18").Value = 1.199056429151248
$ws.Range("D18").Value = 0.09861231489042488
$ws.Range("E18").Value = -0.111662083279245
$ws.Range("F18").Value = 0.4868821201044854
$ws.Range("G18").Value = 0.161213479329271
$ws.Range("H18").Value = 0.8783932807048747
$ws.Range("I18").Value = 2.574265556342576
$ws.Range("J18").Value = 9.418512351267154
$ws.Range("K18").Value = -8.132141731834157

# Row 19
$ws.Range("B19").Value = 1.205589761734299
$ws.Range("C19").Value = 0.1051456474734768
$ws.Range("D19").Value = -0.1051287506961931
$ws.Range("E19").Value = 0.4934154526875373
$ws.Range("F19").Value = 0.1677468119123229
$ws.Range("G19").Value = 0.8849266132879265
$ws.Range("H19").Value = 2.580798888925627
$ws.Range("I19").Value = 9.425045683850206
$ws.Range("J19").Value = -8.125608399251105
$ws.Range("K19").Value = -0.4162117995949584

# Row 20
$ws.Range("B20").Value = -0.006071047505593896
$ws.Range("C20").Value = -0.2163454456752638
$ws.Range("D20").Value = 0.3821987577084666
$ws.Range("E20").Value = 0.05653011693325222
$ws.Range("F20").Value = 0.7737099183088558
$ws.Range("G20").Value = 2.469582193946557
$ws.Range("H20").Value = 9.313828988871135
$ws.Range("I20").Value = -8.236825094230175
$ws.Range("J20").Value = -0.5274284945740291
$ws.Range("K20").Value = 1.055484166312883

# Row 21
$ws.Range("B21").Value = -0.2580915896621678
$ws.Range("C21").Value = 0.3404526137215625
$ws.Range("D21").Value = 0.01478397294634815
$ws.Range("E21").Value = 0.7319637743219518
$ws.Range("F21").Value = 2.427836049959653
$ws.Range("G21").Value = 9.272082844884231
$ws.Range("H21").Value = -8.278571238217079
$ws.Range("I21").Value = -0.5691746385609331
$ws.Range("J21").Value = 1.013738022325978
$ws.Range("K21").Value = -1.952548970023277

# Row 22
$ws.Range("B22").Value = 0.3657676764542774
$ws.Range("C22").Value = 0.04009903567906303
$ws.Range("D22").Value = 0.7572788370546667
$ws.Range("E22").Value = 2.453151112692368
$ws.Range("F22").Value = 9.297397907616945
$ws.Range("G22").Value = -8.253256175484365
$ws.Range("H22").Value = -0.5438595758282182
$ws.Range("I22").Value = 1.039053085058693
$ws.Range("J22").Value = -1.927233907290562
$ws.Range("K22").Value = 0.3367972473739005

# Row 23
$ws.Range("B23").Value = 0.04413770072197692
$ws.Range("C23").Value = 0.7613175020975806
$ws.Range("D23").Value = 2.457189777735282
$ws.Range("E23").Value = 9.30143657265986
$ws.Range("F23").Value = -8.249217510441451
$ws.Range("G23").Value = -0.5398209107853043
$ws.Range("H23").Value = 1.043091750101607
$ws.Range("I23").Value = -1.923195242247649
$ws.Range("J23").Value = 0.3408359124168144
$ws.Range("K23").Value = -0.1103182425099242

# Row 24
$ws.Range("B24").Value = 0.6424403654065582
$ws.Range("C24").Value = 2.338312641044259
$ws.Range("D24").Value = 9.182559435968837
$ws.Range("E24").Value = -8.368094647132473
$ws.Range("F24").Value = -0.6586980474763267
$ws.Range("G24").Value = 0.9242146134105849
$ws.Range("H24").Value = -2.042072378938671
$ws.Range("I24").Value = 0.221958775725792
$ws.Range("J24").Value = -0.2291953792009466
$ws.Range("K24").Value = -0.3695594427149207

# Row 25
$ws.Range("B25").Value = 2.297389002388887
$ws.Range("C25").Value = 9.141635797313464
$ws.Range("D25").Value = -8.409018285787846
$ws.Range("E25").Value = -0.6996216861316987
$ws.Range("F25").Value = 0.8832909747552129
$ws.Range("G25").Value = -2.082996017594043
$ws.Range("H25").Value = 0.18103513707042
$ws.Range("I25").Value = -0.2701190178563186
$ws.Range("J25").Value = -0.4104830813702928
$ws.Range("K25").Value = 0.1157559294919248

# Row 26
$ws.Range("B26").Value = 8.826710628892494
$ws.Range("C26").Value = -8.723943454208817
$ws.Range("D26").Value = -1.01454685455267
$ws.Range("E26").Value = 0.5683658063342414
$ws.Range("F26").Value = -2.397921186015015
$ws.Range("G26").Value = -0.1338900313505515
$ws.Range("H26").Value = -0.5850441862772902
$ws.Range("I26").Value = -0.7254082497912643
$ws.Range("J26").Value = -0.1991692389290468
$ws.Range("K26").Value = -0.3409981856683208

# Row 27
$ws.Range("B27").Value = -9.780318414391347
$ws.Range("C27").Value = -2.0709218147352
$ws.Range("D27").Value = -0.4880091538482882
$ws.Range("E27").Value = -3.454296146197544
$ws.Range("F27").Value = -1.190264991533081
$ws.Range("G27").Value = -1.64141914645982
$ws.Range("H27").Value = -1.781783209973794
$ws.Range("I27").Value = -1.255544199111576
$ws.Range("J27").Value = -1.39737314585085
$ws.Range("K27").Value = -0.7543985244487865

# Row 28
$ws.Range("B28").Value = -1.200275438764269
$ws.Range("C28").Value = 0.3826372221226423
$ws.Range("D28").Value = -2.583649770226613
$ws.Range("E28").Value = -0.3196186155621505
$ws.Range("F28").Value = -0.7707727704888893
$ws.Range("G28").Value = -0.9111368340028634
$ws.Range("H28").Value = -0.3848978231406458
$ws.Range("I28").Value = -0.5267267698799198
$ws.Range("J28").Value = 0.1162478515221441
$ws.Range("K28").Value = 0.2406605578517265

# Row 29
$ws.Range("B29").Value = 0.3719860057927588
$ws.Range("C29").Value = -2.594300986556497
$ws.Range("D29").Value = -0.3302698318920341
$ws.Range("E29").Value = -0.7814239868187727
$ws.Range("F29").Value = -0.9217880503327469
$ws.Range("G29").Value = -0.3955490394705293
$ws.Range("H29").Value = -0.5373779862098034
$ws.Range("I29").Value = 0.1055966351922606
$ws.Range("J29").Value = 0.230009341521843
$ws.Range("K29").Value = -0.1676208645157742

# Row 30
$ws.Range("B30").Value = -2.702915518772638
$ws.Range("C30").Value = -0.4388843641081749
$ws.Range("D30").Value = -0.8900385190349136
$ws.Range("E30").Value = -1.030402582548888
$ws.Range("F30").Value = -0.5041635716866701
$ws.Range("G30").Value = -0.6459925184259441
$ws.Range("H30").Value = -0.003017897023880223
$ws.Range("I30").Value = 0.1213948093057022
$ws.Range("J30").Value = -0.2762353967319151
$ws.Range("K30").Value = -0.3524041104327608

# Row 31
$ws.Range("B31").Value = -0.2307826431404359
$ws.Range("C31").Value = -0.6819367980671746
$ws.Range("D31").Value = -0.8223008615811487
$ws.Range("E31").Value = -0.2960618507189311
$ws.Range("F31").Value = -0.4378907974582051
$ws.Range("G31").Value = 0.2050838239438588
$ws.Range("H31").Value = 0.3294965302734412
$ws.Range("I31").Value = -0.06813367576417605
$ws.Range("J31").Value = -0.1443023894650218
$ws.Range("K31").Value = 0.1320196686578998

# Row 32
$ws.Range("B32").Value = -0.5654386276933741
$ws.Range("C32").Value = -0.7058026912073482
$ws.Range("D32").Value = -0.1795636803451306
$ws.Range("E32").Value = -0.3213926270844047
$ws.Range("F32").Value = 0.3215819943176592
$ws.Range("G32").Value = 0.4459947006472416
$ws.Range("H32").Value = 0.04836449460962439
$ws.Range("I32").Value = -0.02780421909122137
$ws.Range("J32").Value = 0.2485178390317002
$ws.Range("K32").Value = 0.07729801671997623

# Row 33
$ws.Range("B33").Value = -0.6603092772102132
$ws.Range("C33").Value = -0.1340702663479956
$ws.Range("D33").Value = -0.2758992130872696
$ws.Range("E33").Value = 0.3670754083147943
$ws.Range("F33").Value = 0.4914881146443768
$ws.Range("G33").Value = 0.09385790860675949
$ws.Range("H33").Value = 0.01768919490591373
$ws.Range("I33").Value = 0.2940112530288354
$ws.Range("J33").Value = 0.1227914307171113
$ws.Range("K33").Value = 0.2759388362258526

# Row 34
$ws.Range("B34").Value = -0.15162438770796
$ws.Range("C34").Value = -0.293453334447234
$ws.Range("D34").Value = 0.3495212869548299
$ws.Range("E34").Value = 0.4739339932844123
$ws.Range("F34").Value = 0.07630378724679503
$ws.Range("G34").Value = 0.0001350735459492769
$ws.Range("H34").Value = 0.2764571316688709
$ws.Range("I34").Value = 0.1052373093571469
$ws.Range("J34").Value = 0.2583847148658881
$ws.Range("K34").Value = -0.01904405955723064

# Row 35
$ws.Range("B35").Value = -0.2053460154962278
$ws.Range("C35").Value = 0.4376286059058361
$ws.Range("D35").Value = 0.5620413122354185
$ws.Range("E35").Value = 0.1644111061978012
$ws.Range("F35").Value = 0.08824239249695551
$ws.Range("G35").Value = 0.3645644506198771
$ws.Range("H35").Value = 0.1933446283081531
$ws.Range("I35").Value = 0.3464920338168943
$ws.Range("J35").Value = 0.06906325939377558
$ws.Range("K35").Value = -0.09392443396517081

# Row 36
$ws.Range("B36").Value = 0.6162032393936197
$ws.Range("C36").Value = 0.7406159457232021
$ws.Range("D36").Value = 0.3429857396855849
$ws.Range("E36").Value = 0.2668170259847391
$ws.Range("F36").Value = 0.5431390841076607
$ws.Range("G36").Value = 0.3719192617959367
$ws.Range("H36").Value = 0.525066667304678
$ws.Range("I36").Value = 0.2476378928815592
$ws.Range("J36").Value = 0.0846501995226128
$ws.Range("K36").Value = -0.03943237587190501

# Row 37
$ws.Range("B37").Value = 1.652643173475852
$ws.Range("C37").Value = 1.255012967438235
$ws.Range("D37").Value = 1.178844253737389
$ws.Range("E37").Value = 1.455166311860311
$ws.Range("F37").Value = 1.283946489548587
$ws.Range("G37").Value = 1.437093895057328
$ws.Range("H37").Value = 1.159665120634209
$ws.Range("I37").Value = 0.9966774272752628
$ws.Range("J37").Value = 0.8725948518807449
$ws.Range("K37").Value = 1.388747888886706

# Row 38
$ws.Range("B38").Value = 0.3110387314724781
$ws.Range("C38").Value = 0.2348700177716323
$ws.Range("D38").Value = 0.5111920758945538
$ws.Range("E38").Value = 0.3399722535828299
$ws.Range("F38").Value = 0.4931196590915711
$ws.Range("G38").Value = 0.2156908846684524
$ws.Range("H38").Value = 0.05270319130950599
$ws.Range("I38").Value = -0.07137938408501182
$ws.Range("J38").Value = 0.444773652920949
$ws.Range("K38").Value = 0.2348700177716323

# Row 39
$ws.Range("B39").Value = 0.2388379152847414
$ws.Range("C39").Value = 0.5151599734076631
$ws.Range("D39").Value = 0.343940151095939
$ws.Range("E39").Value = 0.4970875566046802
$ws.Range("F39").Value = 0.2196587821815615
$ws.Range("G39").Value = 0.0566710888226151
$ws.Range("H39").Value = -0.06741148657190271
$ws.Range("I39").Value = 0.4487415504340581
$ws.Range("J39").Value = 0.2388379152847414

# Row 40
$ws.Range("B40").Value = 0.6508000635779043
$ws.Range("C40").Value = 0.4795802412661804
$ws.Range("D40").Value = 0.6327276467749217
$ws.Range("E40").Value = 0.3552988723518029
$ws.Range("F40").Value = 0.1923111789928565
$ws.Range("G40").Value = 0.06822860359833866
$ws.Range("H40").Value = 0.5843816406042994
$ws.Range("I40").Value = 0.3744780054549828

# Row 41
$ws.Range("B41").Value = 0.2387740594105157
$ws.Range("C41").Value = 0.3919214649192569
$ws.Range("D41").Value = 0.1144926904961382
$ws.Range("E41").Value = -0.04849500286280822
$ws.Range("F41").Value = -0.172577578257326
$ws.Range("G41").Value = 0.3435754587486348
$ws.Range("H41").Value = 0.1336718235993181

# Row 42
$ws.Range("B42").Value = 0.3465902496671606
$ws.Range("C42").Value = 0.0691614752440418
$ws.Range("D42").Value = -0.09382621811490459
$ws.Range("E42").Value = -0.2179087935094224
$ws.Range("F42").Value = 0.2982442434965384
$ws.Range("G42").Value = 0.08834060834722172

# Row 43
$ws.Range("B43").Value = 0.00230005330798793
$ws.Range("C43").Value = -0.1606876400509585
$ws.Range("D43").Value = -0.2847702154454763
$ws.Range("E43").Value = 0.2313828215604846
$ws.Range("F43").Value = 0.02147918641116785

# Row 44
$ws.Range("B44").Value = -0.1902738424076751
$ws.Range("C44").Value = -0.3143564178021929
$ws.Range("D44").Value = 0.201796619203768
$ws.Range("E44").Value = -0.00810701594554874

# Row 45
$ws.Range("B45").Value = -0.3325070745318338
$ws.Range("C45").Value = 0.1836459624741271
$ws.Range("D45").Value = -0.02625767267518964

# Row 46
$ws.Range("B46").Value = 0.1656141382254278
$ws.Range("C46").Value = -0.04428949692388896

# Row 47
$ws.Range("B47").Value = -0.09587373626955231
